$d = $word.ActiveDocument

$replacements = @(
    @("446÷2=223, 0", "510÷3=170, 0"),
    @("693÷5=138, 3", "338÷7=48, 2"),
    @("765÷4=191, 1", "280÷3=93, 1"),
    @("900÷4=225, 0", "876÷7=125, 1"),
    @("574÷3=191, 1", "543÷9=60, 3"),
    @("268÷7=38, 2", "269÷2=134, 1"),
    @("720÷2=360, 0", "650÷8=81, 2"),
    @("912÷5=182, 2", "629÷5=125, 4"),
    @("566÷5=113, 1", "255÷3=85, 0"),
    @("229÷9=25, 4", "344÷3=114, 2"),
    @("836÷3=278, 2", "231÷4=57, 3"),
    @("597÷2=298, 1", "941÷6=156, 5"),
    @("379÷7=54, 1", "618÷9=68, 6"),
    @("846÷2=423, 0", "497÷4=124, 1"),
    @("846÷6=141, 0", "169÷8=21, 1"),
    @("682÷6=113, 4", "838÷3=279, 1"),
    @("827÷6=137, 5", "605÷8=75, 5"),
    @("960÷3=320, 0", "378÷8=47, 2"),
    @("898÷4=224, 2", "661÷3=220, 1"),
    @("210÷2=105, 0", "499÷3=166, 1"),
    @("825÷4=206, 1", "882÷4=220, 2"),
    @("208÷7=29, 5", "468÷3=156, 0"),
    @("544÷8=68, 0", "171÷8=21, 3"),
    @("344÷6=57, 2", "385÷6=64, 1"),
    @("997÷3=332, 1", "419÷6=69, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
